# Commit: Continuacao desenvolvimento do sample V2
# Implementacao de Permissions e LocalizationText no NavMenu; Implementacao da pagina de Roles.
#
# Adds 16 new sysLocalizationText rows (Logout, AlterProfileImage and the new NavMenu
# entries) to both the "en-us" and "pt-br" worksheets, following the exact same
# formula pattern already used throughout the sheet ($N$1 wraps values in quotes,
# column C mirrors the id, column F builds the final INSERT statement).

$wb = $excel.ActiveWorkbook

# --- shared data for both locales -------------------------------------------------
$keys = @(
    'Logout-Label',`
    'AlterProfileImage-Description',`
    'Menu-Superadmin',`
    'Menu-Configs',`
    'Menu-Instance',`
    'Menu-Roles',`
    'Menu-ObjectPermissions',`
    'Menu-Permissions',`
    'Menu-Users',`
    'Menu-LocalizationTexts',`
    'Menu-GroupParameters',`
    'Menu-Parameters',`
    'Menu-Monitoring',`
    'Menu-SessionLog',`
    'Menu-DataLog',`
    'Menu-ExceptionLog'
)

$enVals = @(
    'Logout',`
    'Drag the file here you wish to upload, or click on buttom above to browse.',`
    'Super Admin',`
    'Business Configs',`
    'Instances',`
    'Roles',`
    'Object Permissions',`
    'Permissions',`
    'Users',`
    'Localization Texts',`
    'Group Parameters',`
    'Parameters',`
    'Monitoring',`
    'Session Log',`
    'Data Log',`
    'Exception Log'
)

$ptVals = @(
    'Sair',`
    'Arraste aqui o arquivo que deseja, ou clique no botão acima para navegar.',`
    'Super Admin',`
    'Configs de Negócio',`
    'Instâncias',`
    'Roles',`
    'Objetos de Permissões',`
    'Permissões',`
    'Usuários',`
    'Textos de Localização',`
    'Grupo de Parâmetros',`
    'Parâmetros',`
    'Monitoramento',`
    'Logs de Acessos',`
    'Logs de Dados',`
    'Logs de Erros'
)

$firstRow = 295
$lastRow  = 310
$enFirstId = 1297
$ptFirstId = 2297

function Fill-LocalizationSheet {
    param($ws, $localeLiteral, $firstId, $vals)

    # Columns A, B, D, E: plain id + locale literal + the two $N$1-wrapped formulas,
    # entered row by row exactly like the rest of the sheet.
    for ($i = 0; $i -lt $keys.Length; $i++) {
        $r = $firstRow + $i
        $ws.Cells.Item($r, 1).Value2 = $firstId + $i
        $ws.Cells.Item($r, 2).Value2 = $localeLiteral
        $ws.Cells.Item($r, 4).Formula = '=$N$1 & "' + $keys[$i] + '" & $N$1'
        $ws.Cells.Item($r, 5).Formula = '=$N$1 & "' + $vals[$i] + '" & $N$1'
    }

    # Rows 295 & 296: the first two new keys (Logout-Label, AlterProfileImage-Description)
    # are typed in one at a time, same as the last few rows already on the sheet.
    $cFormula = '="''" & A295 & "''"'
    $fFormula = '="insert into sysLocalizationText Values(" &A295 & "," & B295 & "," &C295 & "," & D295 & "," & E295 & ",getdate(),getdate())"'
    $ws.Range("C295").Formula = $cFormula
    $ws.Range("F295").Formula = $fFormula

    $cFormula = '="''" & A296 & "''"'
    $fFormula = '="insert into sysLocalizationText Values(" &A296 & "," & B296 & "," &C296 & "," & D296 & "," & E296 & ",getdate(),getdate())"'
    $ws.Range("C296").Formula = $cFormula
    $ws.Range("F296").Formula = $fFormula

    # Rows 297-310 (all the Menu-* entries): filled as a single block, the same way
    # a drag-fill of the formula down the rest of the column would behave.
    $cFormula = '="''" & A297 & "''"'
    $fFormula = '="insert into sysLocalizationText Values(" &A297 & "," & B297 & "," &C297 & "," & D297 & "," & E297 & ",getdate(),getdate())"'
    $ws.Range("C297:C" + $lastRow).Formula = $cFormula
    $ws.Range("F297:F" + $lastRow).Formula = $fFormula
}

$wsEn = $wb.Worksheets.Item("en-us")
$wsPt = $wb.Worksheets.Item("pt-br")

Fill-LocalizationSheet $wsEn "      'en-us'" $enFirstId $enVals
Fill-LocalizationSheet $wsPt "      'pt-br'" $ptFirstId $ptVals

# --- viewport / selection bookkeeping, matching where the author scrolled to next --
$wsEn.Activate()
$wsEn.Range("E303").Select()
$winEn = $excel.ActiveWindow
$winEn.ScrollRow = 286
$winEn.ScrollColumn = 1

$wsPt.Activate()
$wsPt.Range("E301").Select()
$winPt = $excel.ActiveWindow
$winPt.ScrollRow = 289
$winPt.ScrollColumn = 5

$wsEn.Activate()

